$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to reflect new shared strings: Groups, Check, loginTest3
$ws.Range("G1").Value = "Groups"
$ws.Range("G2").Value = "Check"
$ws.Range("F6").Value = "loginTest3"
$ws.Range("G6").Value = "Smoke"

# Update the active selection to G5
$ws.Range("G5").Select()
